# Roll the Jack/Jill retirement-projection templates forward one year
# (2024-start -> 2025-start), dropping the now-superfluous trailing row
# on each sheet and doing a small cleanup on Jack!H3.

$wb = $excel.ActiveWorkbook

# ---- Sheet "Jack" ----------------------------------------------------
$ws1 = $wb.Worksheets.Item("Jack")

# Drop the orphan last row (used to hold only the final year label).
$ws1.Rows.Item(33).Delete()

# The row that is now last (32) used to carry full (empty, styled)
# B:I cells; once it becomes the new "label only" row it should match
# the shape of the row that used to be last.
$ws1.Range("B32:I32").Clear()

# Shift every year in column A forward by one.
for ($r = 2; $r -le 32; $r++) {
    $cell = $ws1.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 + 1
}

# Unrelated small cleanup noted in the diff.
$ws1.Range("H3").Clear()

# ---- Sheet "Jill" ------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Jill")

# Drop the orphan last row (used to hold only the final year label).
$ws2.Rows.Item(36).Delete()

# Row 34 was the last "full" row; once row 35 (now the new last row)
# keeps being label-only, row 34 must also become label-only.
$ws2.Range("B34:I34").Clear()

# Shift every year in column A forward by one.
for ($r = 2; $r -le 35; $r++) {
    $cell = $ws2.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 + 1
}

# ---- Selection / active-sheet bookkeeping ------------------------------
# Jill's selection becomes the whole second row (as if the user had just
# clicked its row header) ...
$ws2.Rows.Item(2).Select() | Out-Null

# ... while Jack becomes the active sheet/tab with E22 selected.
$ws1.Select() | Out-Null
$ws1.Range("E22").Select() | Out-Null
